# claim_details.py: Insert rate_name column and on/off-site summary row
#
# Updates the "sep_2020" monthly column (G) on the crisis_src sheet with the
# counts that were previously recorded as 0, and rolls those same counts
# into the "SFY 2021 Total" column (Q) so the yearly total reflects the
# newly-reported September 2020 data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> sep_2020 (column G) count
$septValues = @{
    3  = 327
    4  = 54
    5  = 381
    8  = 28
    9  = 2
    13 = 15
    14 = 45
    15 = 17
    17 = 100
    18 = 42
    19 = 1
    20 = 2
    21 = 101
    22 = 2
    23 = 2
    24 = 4
    25 = 24
    26 = 42
    27 = 5
    28 = 9
    29 = 34
    31 = 4
    34 = 1
    36 = 9
    38 = 1
    47 = 1
}

# Row 15's yearly total (column Q) was left at its old value in the source
# data, so it is intentionally excluded from the Q-column update below.
$totalValues = $septValues.Clone()
$totalValues.Remove(15)

foreach ($row in $septValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $septValues[$row]
}

foreach ($row in $totalValues.Keys) {
    $ws.Cells.Item($row, 17).Value = $totalValues[$row]
}
